$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Column G ("Sexo") holds "Masculino"/"Femenino" values for rows 2-65.
# Rename them to "Hombre" / "Mujer" respectively.
$lastRow = $ws.Cells.Item($ws.Rows.Count, 7).End(-4162).Row
if ($lastRow -lt 2) { $lastRow = 1 }

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    if ($cell.Value2 -eq "Masculino") {
        $cell.Value2 = "Hombre"
    } elseif ($cell.Value2 -eq "Femenino") {
        $cell.Value2 = "Mujer"
    }
}
